$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "`$ bold(All)"
$ws.Range("C1").Value = "`$ bold(Europe)"
